# Edit script: insert two new weekly records (rows) into the price table.
# The workbook contains a single data sheet where rows 205..310 hold daily/weekly
# price observations for "Palta" (avocado) at "Feria Lagunitas de Puerto Montt".
# This commit adds one new week of observations (date 2022-01-11 / serial 44572)
# for qualities "Primera" and "Segunda", inserted right before the existing
# row 205, pushing all subsequent rows down by two and extending the sheet's
# dimension from A1:T310 to A1:T312.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 205 and 206; everything currently at 205.. shifts down by 2.
$ws.Rows("205:206").Insert()

# --- New row 205: Palta, calidad "Primera" ---
$ws.Range("A205").Value2 = 4
$ws.Range("B205").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C205").Value2 = "Los Lagos"
$ws.Range("D205").Value2 = 44572
$ws.Range("E205").Value2 = 10
$ws.Range("F205").Value2 = "Fruta"
$ws.Range("G205").Value2 = 100106
$ws.Range("H205").Value2 = "Oleaginosos"
$ws.Range("I205").Value2 = 100106002
$ws.Range("J205").Value2 = "Palta"
$ws.Range("K205").Value2 = "Hass"
$ws.Range("L205").Value2 = "Primera"
$ws.Range("M205").Value2 = 400
$ws.Range("N205").Value2 = 4000
$ws.Range("O205").Value2 = 4100
$ws.Range("P205").Value2 = 4050
$ws.Range("Q205").Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R205").Value2 = "Provincia de Quillota"
$ws.Range("S205").Value2 = 4050
$ws.Range("T205").Value2 = 1

# --- New row 206: Palta, calidad "Segunda" ---
$ws.Range("A206").Value2 = 4
$ws.Range("B206").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C206").Value2 = "Los Lagos"
$ws.Range("D206").Value2 = 44572
$ws.Range("E206").Value2 = 10
$ws.Range("F206").Value2 = "Fruta"
$ws.Range("G206").Value2 = 100106
$ws.Range("H206").Value2 = "Oleaginosos"
$ws.Range("I206").Value2 = 100106002
$ws.Range("J206").Value2 = "Palta"
$ws.Range("K206").Value2 = "Hass"
$ws.Range("L206").Value2 = "Segunda"
$ws.Range("M206").Value2 = 200
$ws.Range("N206").Value2 = 3600
$ws.Range("O206").Value2 = 3600
$ws.Range("P206").Value2 = 3600
$ws.Range("Q206").Value2 = "`$/kilo (en caja de 17 kilos)"
$ws.Range("R206").Value2 = "Provincia de Quillota"
$ws.Range("S206").Value2 = 3600
$ws.Range("T206").Value2 = 1
